$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric stat corrections (values were off by one day due to NBA stats timing) ---
$ws.Range("AJ2").Value = 28
$ws.Range("AD4").Value = 8
$ws.Range("AD5").Value = 17
$ws.Range("AD6").Value = 17
$ws.Range("AD7").Value = 17
$ws.Range("AD8").Value = 17
$ws.Range("AH9").Value = 28
$ws.Range("AW9").Value = 17
$ws.Range("AD11").Value = 8
$ws.Range("AQ13").Value = 7
$ws.Range("BB13").Value = 29
$ws.Range("AD15").Value = 8
$ws.Range("AD16").Value = 8
$ws.Range("AD17").Value = 8
$ws.Range("AR17").Value = 15
$ws.Range("AD18").Value = 17
$ws.Range("D19").Value = 39
$ws.Range("F19").Value = 21
$ws.Range("G19").Value = 0.462
$ws.Range("J19").Value = 77.40000000000001
$ws.Range("M19").Value = 16.9
$ws.Range("N19").Value = 0.326
$ws.Range("O19").Value = 20.7
$ws.Range("Q19").Value = 0.718
$ws.Range("R19").Value = 11.8
$ws.Range("S19").Value = 30
$ws.Range("U19").Value = 23.5
$ws.Range("V19").Value = 15.7
$ws.Range("X19").Value = 4.7
$ws.Range("Y19").Value = 4.3
$ws.Range("Z19").Value = 23.2
$ws.Range("AC19").Value = -5
$ws.Range("AD19").Value = 17
$ws.Range("AH19").Value = 5
$ws.Range("AJ19").Value = 27
$ws.Range("AQ19").Value = 24
$ws.Range("AR19").Value = 13
$ws.Range("BB19").Value = 28
$ws.Range("AD20").Value = 17
$ws.Range("AD21").Value = 8
$ws.Range("AQ22").Value = 23
$ws.Range("D24").Value = 40
$ws.Range("E24").Value = 28
$ws.Range("AD24").Value = 8
$ws.Range("AE24").Value = 3
$ws.Range("G24").Value = 0.7
$ws.Range("I24").Value = 41.7
$ws.Range("J24").Value = 85.09999999999999
$ws.Range("L24").Value = 8.9
$ws.Range("M24").Value = 23.3
$ws.Range("N24").Value = 0.384
$ws.Range("O24").Value = 17.5
$ws.Range("P24").Value = 22.5
$ws.Range("Q24").Value = 0.778
$ws.Range("R24").Value = 8.699999999999999
$ws.Range("T24").Value = 41.3
$ws.Range("V24").Value = 13.4
$ws.Range("Y24").Value = 3.8
$ws.Range("AA24").Value = 20.1
$ws.Range("AB24").Value = 109.8
$ws.Range("AC24").Value = 5.4
$ws.Range("AG24").Value = 3
$ws.Range("AN24").Value = 4
$ws.Range("AQ24").Value = 6
$ws.Range("AW24").Value = 18
$ws.Range("AY24").Value = 3
$ws.Range("BC24").Value = 5
$ws.Range("AD25").Value = 8
$ws.Range("AY25").Value = 4
$ws.Range("AD26").Value = 17
$ws.Range("AF26").Value = 19
$ws.Range("AN27").Value = 3
$ws.Range("BC27").Value = 4
$ws.Range("AD28").Value = 8
$ws.Range("AR28").Value = 14

# --- Date column fix: "1-20-2007-08" -> "2008-01-20" (kept as text, not converted to a date serial) ---
$dateCells = @("BF2", "BF3", "BF4", "BF5", "BF6", "BF7", "BF8", "BF9", "BF10", "BF11", "BF12", "BF13", "BF14", "BF15", "BF16", "BF17", "BF18", "BF19", "BF20", "BF21", "BF22", "BF23", "BF24", "BF25", "BF26", "BF27", "BF28", "BF29", "BF30", "BF31")
foreach ($cell in $dateCells) {
    $origStyle = $ws.Range($cell).Style
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = "2008-01-20"
    $ws.Range($cell).Style = $origStyle
}
